# Apply the two changes described by the commit:
#  1) Re-point every table's style (currently the GUID
#     {26083FAA-C9F0-4969-8330-AAEF655F36DE}) to the new table style GUID
#     {EDDECFFA-6313-4EA5-B4C0-B2114FA61464}.
#  2) Swap the "Integral" (Red Violet) / "Office Theme" colour palettes
#     that live in the deck's two theme parts, i.e. the slide master's
#     theme takes on the Office Theme palette that used to belong to the
#     notes master's theme.

$p = $ppt.ActivePresentation

# --- 1) Table styles -------------------------------------------------
$oldStyle = "{26083FAA-C9F0-4969-8330-AAEF655F36DE}"
$newStyle = "{EDDECFFA-6313-4EA5-B4C0-B2114FA61464}"

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shp = $slide.Shapes.Item($j)
        if ($shp.HasTable) {
            $tbl = $shp.Table
            if ($tbl.Style.Name -eq $oldStyle) {
                $tbl.ApplyStyle($newStyle)
            }
        }
    }
}

# --- 2) Theme colour swap --------------------------------------------
# The presentation's theme colour scheme (12 slots: dk1, lt1, dk2, lt2,
# accent1-6, hlink, folHlink) moves from the "Integral" / Red Violet
# palette to the "Office Theme" palette (the palette previously used by
# the notes master's theme part).
$officePalette = @(
    0x000000,  # dk1
    0xFFFFFF,  # lt1
    0x44546A,  # dk2
    0xE7E6E6,  # lt2
    0x5B9BD5,  # accent1
    0xED7D31,  # accent2
    0xA5A5A5,  # accent3
    0xFFC000,  # accent4
    0x4472C4,  # accent5
    0x70AD47,  # accent6
    0x0563C1,  # hlink
    0x954F72   # folHlink
)

$slide1 = $p.Slides.Item(1)
$themeColors = $slide1.ThemeColorScheme

for ($i = 1; $i -le $themeColors.Count; $i++) {
    $rgbHex = $officePalette[$i - 1]
    $r = [math]::Floor($rgbHex / 0x10000) -band 0xFF
    $g = [math]::Floor($rgbHex / 0x100) -band 0xFF
    $b = $rgbHex -band 0xFF
    # The RGB property is packed as 0xBBGGRR (classic OLE RGB()).
    $packed = ($b * 0x10000) + ($g * 0x100) + $r
    $themeColors.Colors($i).RGB = $packed
}
